# Add info for REV051 and REV138 to systematicWrongButtons.
# These participants were using 91 for right and 94 for left and so had
# accumulated a lot of Failed Gos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows right after the existing "51" block (old row 40)
# so the two new rows become rows 41-42, matching the target layout. ---
$ws.Rows("41:42").Insert()

# Row 41: subject 51, trial 1, left=94, right=91
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "51"
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = 94
$ws.Cells.Item(41, 4).Value = 91

# Row 42: subject 51, trial 2, left=94, right=91
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "51"
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = 94
$ws.Cells.Item(42, 4).Value = 91

# --- Append two new rows at the bottom of the table for subject 138 ---
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "138"
$ws.Cells.Item(71, 2).Value = 1
$ws.Cells.Item(71, 3).Value = 94
$ws.Cells.Item(71, 4).Value = 91

$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "138"
$ws.Cells.Item(72, 2).Value = 2
$ws.Cells.Item(72, 3).Value = 94
$ws.Cells.Item(72, 4).Value = 91

# --- Match the final view state recorded in the saved workbook ---
$ws.Range("C73").Select()
